$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tool Settings")

# Insert a new row above the current row 8 ("{{Inference_Tool:...}}" row),
# shifting the rows below it down by one.
$ws.Rows("7").Insert()

# New explanatory note about re-using a tool with multiple settings.
$noteText = "If the same tool has been used with multiple settings (e.g. with different starting databases), please create a table for each setting used. " + [char]10 + "Please use a unique tool_name for each tool and setting combination, and use this name when referring to the tool in the Genotype table."
$ws.Range("B7").Value = $noteText
$ws.Range("B7").WrapText = $false
$ws.Rows("7").AutoFit()

# Make "Tool Settings" the active sheet/tab, matching the saved selection.
$ws.Activate() | Out-Null
$ws.Range("G15").Select() | Out-Null
